# Regenerate the saved pitching-log data for column G ("K" - strikeouts).
# The source data previously stored a different stat ("Strike#") in column
# G; it is being regenerated to hold the true strikeout count (K) per
# outing, after recalculating the underlying std/mean stats upstream and
# writing the corrected s_vals back into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (strikeouts) value, taken from the
# regenerated save data.
$kValues = [ordered]@{
    2  = 0
    3  = 1
    4  = 1
    5  = 2
    6  = 1
    7  = 0
    8  = 1
    9  = 0
    10 = 2
    11 = 3
    12 = 0
    13 = 1
    14 = 0
    15 = 2
    16 = 0
    17 = 1
    18 = 0
    19 = 2
    20 = 0
    21 = 2
    22 = 1
    23 = 1
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 2
    29 = 1
    30 = 0
    31 = 4
    32 = 2
    33 = 0
    34 = 0
    35 = 1
    36 = 1
    37 = 2
    38 = 2
    39 = 0
    40 = 0
    41 = 0
    42 = 1
    43 = 2
    44 = 1
    46 = 1
    47 = 3
    48 = 0
    49 = 0
    50 = 0
    51 = 1
    52 = 1
    53 = 1
    54 = 0
    55 = 1
    56 = 1
    57 = 2
    58 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
